$d = $word.ActiveDocument
$marker = "<1T>"

# ------------------------------------------------------------------
# 1. Build a brand-new, formatting-free paragraph containing the
#    marker text. Paragraphs.Add() (no args) appends a clean
#    <w:p><w:r></w:r></w:p> with no inherited pPr/rPr, unlike
#    InsertParagraphBefore() which clones neighbouring formatting.
#    We give it one extra placeholder character ("X") so that the
#    bookmark we add below lands at a genuine "middle of paragraph"
#    offset instead of a paragraph boundary (boundary offsets get
#    snapped to "whole paragraph" by this bookmark implementation).
# ------------------------------------------------------------------
$d.Paragraphs.Add() | Out-Null
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = $marker + "X"

# ------------------------------------------------------------------
# 2. Move that paragraph (including its paragraph mark) from the end
#    of the document to the very beginning via Cut/Paste.
# ------------------------------------------------------------------
$movRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$movRange.Cut()
$d.Range(0, 0).Paste()

# ------------------------------------------------------------------
# 3. Remove the old "_GoBack" bookmark from the paragraph at the end
#    of the document (it becomes a bare empty paragraph).
# ------------------------------------------------------------------
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# ------------------------------------------------------------------
# 4. Re-create the "_GoBack" bookmark, collapsed, right after the
#    marker text (before the trailing placeholder "X") in the new
#    first paragraph.
# ------------------------------------------------------------------
$bmPoint = $marker.Length
$bmRange = $d.Range($bmPoint, $bmPoint)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ------------------------------------------------------------------
# 5. Drop the placeholder "X" now that the bookmark is anchored.
# ------------------------------------------------------------------
$placeholder = $d.Range($bmPoint, $bmPoint + 1)
$placeholder.Delete()
